# Update the El Salvador MSME summary figures with refined (2-decimal)
# values while keeping the cells as text (matching the original workbook,
# where these figures are stored as shared-string text, not numbers).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$updates = @{
    "B11" = "24.89"   # Enterprises density (per 1000 people) - Micro
    "C11" = "0.92"    # Enterprises density (per 1000 people) - SMEs
    "B12" = "50.81"   # Employment (% of total) - Micro
    "C12" = "23.21"   # Employment (% of total) - SMEs
    "D12" = "74.02"   # Employment (% of total) - MSMEs
    "B14" = "96.16"   # Enterprises (% of total) - Micro
    "C14" = "3.54"    # Enterprises (% of total) - SMEs
    "D14" = "99.69"   # Enterprises (% of total) - MSMEs
}

# Force a Text number format first so the numeric-looking strings are
# written back as literal text (shared strings) instead of being
# auto-coerced into numeric cells, then clear the temporary formatting so
# the cell style reverts to the original (default) style.
foreach ($addr in $updates.Keys) {
    $ws.Range($addr).NumberFormat = "@"
}
foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
foreach ($addr in $updates.Keys) {
    $ws.Range($addr).ClearFormats()
}
